$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.806.10'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.529.92'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '624.43'
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").Value = '174.40'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.612'
$ws.Range("E7").Value = '  -0.60%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.526.53'
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("D11").Value = '7.12'
$ws.Range("E11").Value = '  -4.48%  '
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").Value = '46.54'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").Value = '0.0000277'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '4.095.03'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").Value = '8.42'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = '608.03'
$ws.Range("E17").Value = '  -1.40%  '
$ws.Range("D18").Value = '3.528.98'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '70.814.18'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '0.122'
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").Value = '17.78'
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("D22").Value = '0.886'
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").Value = '9.12'
$ws.Range("E23").Value = '  -2.75%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '15.68'
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '98.40'
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("E26").Value = '  -0.92%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '2.58'
$ws.Range("E28").Value = '  -2.40%  '
$ws.Range("D29").Value = '33.84'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '9.09'
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("D31").Value = '3.04'
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -4.03%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").Value = '634.98'
$ws.Range("E34").Value = '  +3.56%  '
$ws.Range("D35").Value = '6.81'
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.0999'
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("B37").Value = 'Cosmos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D37").Value = '10.84'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  -8.46%  '
$ws.Range("D39").Value = '0.0476'
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("D40").Value = '56.83'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("D43").Value = '3.364.63'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").Value = '0.0₃0733'
$ws.Range("E44").Value = '  +2.98%  '
$ws.Range("D45").Value = '2.99'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").Value = '0.312'
$ws.Range("E46").Value = '  -2.76%  '
$ws.Range("D47").Value = '32.09'
$ws.Range("E47").Value = '  -2.82%  '
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").Value = '132.95'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = '0.157'
$ws.Range("E51").Value = '  +6.11%  '
